$d = $word.ActiveDocument

# Locate the target paragraph by finding its distinctive original sentence.
# (Built with an escaped "e with circumflex" to dodge any console/codepage
# round-tripping issues for non-ASCII source characters.)
$ea = [char]0xEA
$needle = "pelos tr" + $ea + "s, ou ainda por nenhum."

$rng = $d.Content
$f = $rng.Find
$f.ClearFormatting()
$f.Text = $needle
$f.MatchWildcards = $false
$found = $f.Execute()
if (-not $found) {
    throw "Could not locate target sentence in document."
}

$para = $rng.Paragraphs(1)
$prng = $para.Range

# Pull this paragraph's real OOXML (preserves w:pPr / paragraph identity
# attributes exactly) rather than hand-authoring them.
$wholeXml = $prng.WordOpenXML
if ($wholeXml -notmatch '(?s)(<w:p [^>]*>.*?</w:p>)') {
    throw "Could not extract paragraph OOXML."
}
$paraXml = $matches[1]

# Build the original single run's full text out of char codes so this file
# stays plain-ASCII while still matching the accented source text exactly.
$a = [char]0xE1   # a with acute
$u = [char]0xFA   # u with acute
$e1 = [char]0xE9  # e with acute (1st set)
$i1 = [char]0xED  # i with acute (1st set)
$atil = [char]0xE3 # a with tilde
$oldText = "Crie um formul" + $a + "rio onde seja solicitado um n" + $u + "mero, atrav" + $e1 + "s de PHP e verifique se o valor " + $e1 + " divis" + $i1 + "vel por 10, por 5, por 2 ou se n" + $atil + "o " + $e1 + " divis" + $i1 + "vel por nenhum deles, retornando se foi divis" + $i1 + "vel por 10 ou por 5 ou por 2 ou pelos tr" + $ea + "s, ou ainda por nenhum."

$oldRun = "<w:r><w:t>" + $oldText + "</w:t></w:r>"
if ($paraXml.IndexOf($oldRun) -lt 0) {
    throw "Could not find the expected original run inside paragraph OOXML."
}

# Split that single run into five runs per the target revision.
$newRuns = "<w:r><w:t>Crie um formul" + $a + "rio onde seja solicitado um n" + $u + "mero, atrav" + $e1 + "s de PHP e verifique se o valor " + $e1 + " divis" + $i1 + "vel por 10, por 5, por 2 ou se n" + $atil + "o " + $e1 + " divis" + $i1 + "vel por nenhum deles, retornando se foi divis" + $i1 + "vel por 10 ou por 5 ou por 2 ou pelos tr" + $ea + "s,</w:t></w:r>" `
    + '<w:r><w:t xml:space="preserve"> r</w:t></w:r>' `
    + "<w:r><w:t>esto da divis" + $atil + "o 0</w:t></w:r>" `
    + '<w:r><w:t>,</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> ou ainda por nenhum.</w:t></w:r>'

$newParaXml = $paraXml.Replace($oldRun, $newRuns)

# Replace the whole paragraph's contents (this op requires the entire
# paragraph range - narrower ranges corrupt neighboring content).
$null = $prng.InsertXML($newParaXml)
